$wb = $excel.ActiveWorkbook

# Map of row -> new "想去人数" (F column) value, identical update applied
# to both the "展览" (sheet1) and "全部类型" (sheet4) worksheets.
$updates = @{
    2  = 138
    3  = 416
    4  = 12031
    5  = 1254
    6  = 127
    10 = 180
    13 = 62
    17 = 2036
    19 = 925
    20 = 119
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
